# Fix a stray trailing apostrophe in the KBS xpath string (B16) and add a
# new blank row (18) below the existing table, carrying the same
# (quote-prefixed) formatting that B16 had.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Clone B16's current formatting (quote-prefix style) down onto the new
#    B18 cell before we touch B16's value.
$ws.Range("B16").Copy()
$ws.Range("B18").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# 2) Correct the KBS xpath value: drop the trailing "'" typo.
$ws.Range("B16").Value = '//*[@id="cont_newstext"]'

# 3) Re-apply the quote-prefix formatting to B16 (writing .Value above reset
#    it to the default style), using B18 - which now carries that format -
#    as the source.
$ws.Range("B18").Copy()
$ws.Range("B16").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# 4) Move the active selection to B17, matching the post-edit sheet view.
$ws.Range("B17").Select()
